# Auto-generated Excel COM-interop script to apply numeric corrections
# to the "Leve Profits" worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each block below updates specific H..N cells for a given row, matching
# refreshed market-board price data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2857
$ws.Range("J32").Value = 2947
$ws.Range("L32").Value = 2947
$ws.Range("N32").Value = -3599

$ws.Range("H40").Value = 1706.6207
$ws.Range("I40").Value = 1747.4783
$ws.Range("J40").Value = 1550
$ws.Range("K40").Value = 1747.4783
$ws.Range("L40").Value = 1550
$ws.Range("M40").Value = -1572.4783
$ws.Range("N40").Value = -1900

$ws.Range("H98").Value = 893.1429000000001
$ws.Range("I98").Value = 899.6667
$ws.Range("K98").Value = 899.6667
$ws.Range("M98").Value = 598.3333

$ws.Range("H100").Value = 1460.2
$ws.Range("I100").Value = 1500
$ws.Range("J100").Value = 1450.25
$ws.Range("K100").Value = 1500
$ws.Range("L100").Value = 1450.25
$ws.Range("M100").Value = -959
$ws.Range("N100").Value = -2532.25

$ws.Range("H103").Value = 1284
$ws.Range("J103").Value = 1341
$ws.Range("L103").Value = 4023
$ws.Range("N103").Value = -5195

$ws.Range("H122").Value = 893.1429000000001
$ws.Range("I122").Value = 899.6667
$ws.Range("K122").Value = 2699.0001
$ws.Range("M122").Value = -249.0001000000002

$ws.Range("H138").Value = 3998.205
$ws.Range("J138").Value = 4446.107
$ws.Range("L138").Value = 13338.321
$ws.Range("N138").Value = -23618.321

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4012.1538
$ws.Range("I45").Value = 4669
$ws.Range("J45").Value = 399.5
$ws.Range("K45").Value = 4669
$ws.Range("L45").Value = 399.5
$ws.Range("M45").Value = -4292
$ws.Range("N45").Value = -1153.5

$ws.Range("H61").Value = 8019.9287
$ws.Range("I61").Value = 8598.25
$ws.Range("K61").Value = 8598.25
$ws.Range("M61").Value = -8386.25

$ws.Range("H74").Value = 1068.091
$ws.Range("J74").Value = 1250
$ws.Range("L74").Value = 1250
$ws.Range("N74").Value = -2998

$ws.Range("H77").Value = 1068.091
$ws.Range("J77").Value = 1250
$ws.Range("L77").Value = 6250
$ws.Range("N77").Value = -14986

$ws.Range("H97").Value = 1721
$ws.Range("I97").Value = 1861.3334
$ws.Range("J97").Value = 1300
$ws.Range("K97").Value = 1861.3334
$ws.Range("L97").Value = 1300
$ws.Range("M97").Value = -1365.3334
$ws.Range("N97").Value = -2292

$ws.Range("H136").Value = 8019.9287
$ws.Range("I136").Value = 8598.25
$ws.Range("K136").Value = 25794.75
$ws.Range("M136").Value = -23244.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 169.16667
$ws.Range("J80").Value = 163.66667
$ws.Range("L80").Value = 163.66667
$ws.Range("N80").Value = -2159.66667

$ws.Range("H83").Value = 169.16667
$ws.Range("J83").Value = 163.66667
$ws.Range("L83").Value = 818.3333500000001
$ws.Range("N83").Value = -10802.33335

$ws.Range("H86").Value = 8363.625
$ws.Range("I86").Value = 8802
$ws.Range("J86").Value = 7633
$ws.Range("K86").Value = 8802
$ws.Range("L86").Value = 7633
$ws.Range("M86").Value = -7679
$ws.Range("N86").Value = -9879

$ws.Range("H89").Value = 8363.625
$ws.Range("I89").Value = 8802
$ws.Range("J89").Value = 7633
$ws.Range("K89").Value = 44010
$ws.Range("L89").Value = 38165
$ws.Range("M89").Value = -38394
$ws.Range("N89").Value = -49397

$ws.Range("H94").Value = 2397.6667
$ws.Range("I94").Value = 2370
$ws.Range("J94").Value = 2494.5
$ws.Range("K94").Value = 2370
$ws.Range("L94").Value = 2494.5
$ws.Range("M94").Value = -1919
$ws.Range("N94").Value = -3396.5

$ws.Range("H105").Value = 2452.75
$ws.Range("I105").Value = 2133.3333
$ws.Range("K105").Value = 2133.3333
$ws.Range("M105").Value = -386.3332999999998

$ws.Range("H107").Value = 7915
$ws.Range("I107").Value = 7372.5
$ws.Range("J107").Value = 9000
$ws.Range("K107").Value = 7372.5
$ws.Range("L107").Value = 9000
$ws.Range("M107").Value = -5452.5
$ws.Range("N107").Value = -12840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6335.5557
$ws.Range("I16").Value = 5503.3335
$ws.Range("J16").Value = 8000
$ws.Range("K16").Value = 5503.3335
$ws.Range("L16").Value = 8000
$ws.Range("M16").Value = -5216.3335
$ws.Range("N16").Value = -8574

$ws.Range("H31").Value = 1225.08
$ws.Range("I31").Value = 1097.8572
$ws.Range("K31").Value = 1097.8572
$ws.Range("M31").Value = -802.8571999999999

$ws.Range("H34").Value = 1225.08
$ws.Range("I34").Value = 1097.8572
$ws.Range("K34").Value = 1097.8572
$ws.Range("M34").Value = -895.8571999999999

$ws.Range("H58").Value = 4627.1665
$ws.Range("I58").Value = 3052.25
$ws.Range("K58").Value = 3052.25
$ws.Range("M58").Value = -2849.25

$ws.Range("H113").Value = 6335.5557
$ws.Range("I113").Value = 5503.3335
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 5503.3335
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = -3333.3335
$ws.Range("N113").Value = -12340

$ws.Range("H136").Value = 4627.1665
$ws.Range("I136").Value = 3052.25
$ws.Range("K136").Value = 9156.75
$ws.Range("M136").Value = -6606.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 5798.5713
$ws.Range("I6").Value = 8113.4
$ws.Range("J6").Value = 11.5
$ws.Range("K6").Value = 24340.2
$ws.Range("L6").Value = 34.5
$ws.Range("M6").Value = -24227.2
$ws.Range("N6").Value = -260.5

$ws.Range("H12").Value = 51
$ws.Range("J12").Value = 57
$ws.Range("L12").Value = 171
$ws.Range("N12").Value = -517

$ws.Range("H125").Value = 3998.75
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H131").Value = 2558.1765
$ws.Range("J131").Value = 2699.5334
$ws.Range("L131").Value = 8098.600199999999
$ws.Range("N131").Value = -18178.6002

$ws.Range("H133").Value = 15998.3
$ws.Range("I133").Value = 6996.6
$ws.Range("J133").Value = 25000
$ws.Range("K133").Value = 20989.8
$ws.Range("L133").Value = 75000
$ws.Range("M133").Value = -15929.8
$ws.Range("N133").Value = -85120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H113").Value = 1303.3334
$ws.Range("I113").Value = 1272.5
$ws.Range("J113").Value = 1365
$ws.Range("K113").Value = 1272.5
$ws.Range("L113").Value = 1365
$ws.Range("M113").Value = 897.5
$ws.Range("N113").Value = -5705

$ws.Range("H126").Value = 5719.6665
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 5719.6665
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 17158.9995
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -22098.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2339.111
$ws.Range("I40").Value = 1840.909
$ws.Range("J40").Value = 3122
$ws.Range("K40").Value = 1840.909
$ws.Range("L40").Value = 3122
$ws.Range("M40").Value = -1704.909
$ws.Range("N40").Value = -3394

$ws.Range("H61").Value = 2204.6
$ws.Range("I61").Value = 2157.1428
$ws.Range("K61").Value = 2157.1428
$ws.Range("M61").Value = -1955.1428

$ws.Range("H82").Value = 977.1429000000001
$ws.Range("I82").Value = 938.4
$ws.Range("J82").Value = 1074
$ws.Range("K82").Value = 938.4
$ws.Range("L82").Value = 1074
$ws.Range("M82").Value = -577.4
$ws.Range("N82").Value = -1796

$ws.Range("H85").Value = 977.1429000000001
$ws.Range("I85").Value = 938.4
$ws.Range("J85").Value = 1074
$ws.Range("K85").Value = 938.4
$ws.Range("L85").Value = 1074
$ws.Range("M85").Value = 309.6
$ws.Range("N85").Value = -3570

$ws.Range("H93").Value = 2699.1667
$ws.Range("I93").Value = 2659.2
$ws.Range("J93").Value = 2899
$ws.Range("K93").Value = 2659.2
$ws.Range("L93").Value = 2899
$ws.Range("M93").Value = -1411.2
$ws.Range("N93").Value = -5395

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H113").Value = 2204.6
$ws.Range("I113").Value = 2157.1428
$ws.Range("K113").Value = 2157.1428
$ws.Range("M113").Value = 12.85719999999992

$ws.Range("H136").Value = 2450.3
$ws.Range("I136").Value = 2324.75
$ws.Range("J136").Value = 2534
$ws.Range("K136").Value = 6974.25
$ws.Range("L136").Value = 7602
$ws.Range("M136").Value = -4424.25
$ws.Range("N136").Value = -12702

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 16390.834
$ws.Range("J4").Value = 22686.25
$ws.Range("L4").Value = 22686.25
$ws.Range("N4").Value = -22912.25

$ws.Range("H122").Value = 4283.2104
$ws.Range("I122").Value = 3015.8462
$ws.Range("K122").Value = 9047.5386
$ws.Range("M122").Value = -6597.5386

$ws.Range("H136").Value = 3781.3157
$ws.Range("I136").Value = 3865.8462
$ws.Range("K136").Value = 11597.5386
$ws.Range("M136").Value = -9047.5386
